# Updates the cryptos price/volume table (columns D and E) with refreshed values.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($range, $text) {
    # Leading apostrophe forces Excel to store the value as text, just like
    # manually typing it in, so purely numeric-looking strings (e.g. "522.27")
    # are not reinterpreted as numbers.
    $range.Value = "'" + $text
    $range.Style = 'Normal'
}

$ws.Range('D2').Value = '59.017.19'
$ws.Range('E2').Value = '  +2.94%  '
$ws.Range('D3').Value = '3.108.64'
$ws.Range('E4').Value = '  +0.01%  '
Set-TextValue $ws.Range('D5') '522.27'
$ws.Range('E5').Value = '  +1.60%  '
Set-TextValue $ws.Range('D6') '144.17'
$ws.Range('E6').Value = '  +2.15%  '
$ws.Range('E7').Value = '  +0.02%  '
$ws.Range('E8').Value = '  +1.01%  '
Set-TextValue $ws.Range('D9') '7.39'
$ws.Range('E9').Value = '  +1.85%  '
$ws.Range('E10').Value = '  +1.09%  '
$ws.Range('E11').Value = '  +3.04%  '
$ws.Range('D12').Value = '3.642.04'
$ws.Range('E12').Value = '  +1.43%  '
$ws.Range('E13').Value = '  +1.40%  '
Set-TextValue $ws.Range('D14') '27.16'
$ws.Range('E14').Value = '  +6.57%  '
$ws.Range('E15').Value = '  +1.42%  '
$ws.Range('D16').Value = '58.970.37'
$ws.Range('E16').Value = '  +2.80%  '
$ws.Range('D17').Value = '3.110.95'
$ws.Range('E17').Value = '  +1.45%  '
Set-TextValue $ws.Range('D18') '6.22'
$ws.Range('E18').Value = '  +3.36%  '
Set-TextValue $ws.Range('D19') '13.06'
$ws.Range('E19').Value = '  +0.57%  '
Set-TextValue $ws.Range('D20') '8.24'
$ws.Range('E20').Value = '  +1.51%  '
Set-TextValue $ws.Range('D21') '343.74'
$ws.Range('E21').Value = '  +1.79%  '
$ws.Range('E22').Value = '  -0.08%  '
Set-TextValue $ws.Range('D23') '0.509'
$ws.Range('E23').Value = '  +2.06%  '
Set-TextValue $ws.Range('D24') '65.76'
$ws.Range('E24').Value = '  +0.49%  '
$ws.Range('E25').Value = '  +0.56%  '
$ws.Range('E26').Value = '  -0.06%  '
$ws.Range('D27').Value = '0.0₃0935'
$ws.Range('E27').Value = '  -0.14%  '
$ws.Range('E28').Value = '  +4.70%  '
Set-TextValue $ws.Range('D29') '7.28'
$ws.Range('E29').Value = '  +2.65%  '
$ws.Range('E30').Value = '  +2.13%  '
Set-TextValue $ws.Range('D31') '1.21'
$ws.Range('E31').Value = '  +3.27%  '
Set-TextValue $ws.Range('D32') '21.05'
$ws.Range('E32').Value = '  +1.47%  '
Set-TextValue $ws.Range('D33') '155.05'
$ws.Range('E33').Value = '  +0.40%  '
Set-TextValue $ws.Range('D34') '4.66'
$ws.Range('E34').Value = '  +3.23%  '
Set-TextValue $ws.Range('D35') '6.21'
$ws.Range('E35').Value = '  +5.85%  '
Set-TextValue $ws.Range('D36') '26.96'
$ws.Range('E36').Value = '  +4.20%  '
$ws.Range('E37').Value = '  +6.18%  '
Set-TextValue $ws.Range('D38') '0.0687'
$ws.Range('E38').Value = '  +1.98%  '
Set-TextValue $ws.Range('D39') '3.96'
$ws.Range('E39').Value = '  +2.77%  '
$ws.Range('D40').Value = '3.150.52'
$ws.Range('E40').Value = '  +1.56%  '
Set-TextValue $ws.Range('D41') '36.93'
$ws.Range('E41').Value = '  -0.23%  '
$ws.Range('E42').Value = '  -0.04%  '
Set-TextValue $ws.Range('D43') '0.667'
$ws.Range('E43').Value = '  -0.25%  '
Set-TextValue $ws.Range('D44') '1.46'
$ws.Range('E44').Value = '  +5.61%  '
$ws.Range('D45').Value = '2.284.50'
$ws.Range('E45').Value = '  +1.46%  '
$ws.Range('E46').Value = '  +2.42%  '
Set-TextValue $ws.Range('D47') '21.01'
$ws.Range('E47').Value = '  +4.53%  '
Set-TextValue $ws.Range('D48') '0.962'
$ws.Range('E48').Value = '  +1.16%  '
$ws.Range('E49').Value = '  +3.08%  '
Set-TextValue $ws.Range('D50') '0.756'
$ws.Range('E50').Value = '  +10.99%  '
Set-TextValue $ws.Range('D51') '263.04'
$ws.Range('E51').Value = '  +11.35%  '
